$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Cargo" column header (bold, matching the other header cells)
$ws.Range("F1").Value = "Cargo"
$ws.Range("F1").Font.Bold = $true

# National candidates (rows 2-7) run for "Presidente"
$ws.Range("F2:F7").Value = "Presidente"

# Provincial candidates (rows 8-23) run for "Gobernador"
$ws.Range("F8:F23").Value = "Gobernador"

# Update the view: scroll so row 5 is at top, and select F8:F23 (active cell F8)
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("F8:F23").Select()
